$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column O header cell (empty, but formatted) - extends used range to O1
$ws.Cells.Item(1, 15).NumberFormat = "@"

# Give the new data rows (102/103) a consistent text format before filling them in,
# so the written values land on a single shared style.
$ws.Range("A102:F103").NumberFormat = "@"
$ws.Range("L102:M103").NumberFormat = "@"

# Row 102 - E. Vatcher (elder)
# Row 103 - S. Vatcher (sister)
# Write in the same order the original author typed them so the shared-string
# table grows in the same sequence.
$ws.Cells.Item(102, 2).Value = "0905388991"
$ws.Cells.Item(103, 2).Value = "0905388992"

$ws.Cells.Item(102, 4).Value = "屈長老"
$ws.Cells.Item(103, 4).Value = "屈姐妹"

$ws.Cells.Item(102, 5).Value = "DAAN_FUXING_COUPLE"
$ws.Cells.Item(103, 5).Value = "DAAN_FUXING_COUPLE"

$ws.Cells.Item(102, 3).Value = "E. Vatcher"
$ws.Cells.Item(102, 6).Value = "E. Vatcher"
$ws.Cells.Item(103, 3).Value = "S. Vatcher"
$ws.Cells.Item(103, 6).Value = "S. Vatcher"

$ws.Cells.Item(102, 1).Value = "VATCHER_E"
$ws.Cells.Item(103, 1).Value = "VATCHER_S"

$ws.Cells.Item(102, 12).Value = "0"
$ws.Cells.Item(102, 13).Value = "0"
$ws.Cells.Item(103, 12).Value = "0"
$ws.Cells.Item(103, 13).Value = "0"

# Match the author's final on-screen selection.
$ws.Range("L89").Select()
